# Applies updated "Price" (D) and "Volume(1h)" (E) values to the cryptos
# listing sheet, as produced by the scheduled GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.924.83"
$ws.Range("E2").Value = "  +5.79%  "
$ws.Range("D3").Value = "2.234.60"
$ws.Range("E3").Value = "  +3.01%  "
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E4").Value = "  +0.12%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "231.64"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +2.19%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "0.626"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +0.58%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "61.79"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -2.09%  "
$ws.Range("E8").Value = "  +0.13%  "
$ws.Range("E9").Value = "  +2.97%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "59.27"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +1.40%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.0893"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +4.81%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.104"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -0.03%  "
$ws.Range("D13").Value = "2.566.03"
$ws.Range("E13").Value = "  +3.06%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "15.67"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -1.37%  "
$ws.Range("E15").Value = "  +1.33%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "0.803"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -0.80%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "5.59"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +1.92%  "
$ws.Range("D18").Value = "2.247.62"
$ws.Range("E18").Value = "  +3.67%  "
$ws.Range("D19").Value = "41.884.91"
$ws.Range("E19").Value = "  +5.80%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "72.12"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +0.53%  "
$ws.Range("D21").Value = "0.0₃0898"
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "6.06"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +0.70%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "249.96"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +9.75%  "
$ws.Range("E24").Value = "  +0.02%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "2.40"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +2.45%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "2.37"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +1.87%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "9.69"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +2.22%  "
$ws.Range("E28").Value = "  +2.67%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "166.70"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -2.43%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "19.97"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -1.35%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "2.63"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -1.56%  "
$ws.Range("E33").Value = "  -0.08%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "5.01"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +6.37%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "4.69"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +3.90%  "
$ws.Range("E36").Value = "  +3.26%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "6.67"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -4.42%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "3.66"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -3.74%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "2.37"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -0.94%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.000256"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +31.07%  "
$ws.Range("E41").Value = "  +0.19%  "
$ws.Range("E42").Value = "  +5.22%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "4.85"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -2.51%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "8.57"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +8.78%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.0978"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +6.16%  "
$ws.Range("E46").Value = "  +0.98%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "98.91"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -3.49%  "
$ws.Range("D48").Value = "1.478.82"
$ws.Range("E48").Value = "  -2.36%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "16.55"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -6.36%  "
$ws.Range("E50").Value = "  +0.29%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "52.55"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +8.80%  "
